$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13.96128838361976, 6.876005029580077, 0, 14.13411350873235, 41.93243388613818, 44.1308114787182, 18.21598862651263, 28.93243347951861, 8.90964630204717, 11.06736725554729, 12.20765481881889, 16.2649007694563, 21.77714851554531),
    @(13.79038372623333, 6.80679284274449, 0, 14.14622629166904, 41.95994156897194, 44.2070400650726, 18.26383573636536, 29.01419969977591, 8.90752245462906, 10.94943087262273, 12.20954885147832, 16.24090924337228, 21.83983895019549),
    @(13.68726588670816, 6.763102746820334, 0, 14.15514013341997, 41.98521945794436, 44.26426561438453, 18.29581873911221, 29.06897842832023, 8.90624671778926, 10.87861097243451, 12.21223158305875, 16.22852299349589, 21.88014434602323),
    @(13.64575048869756, 6.745004312178082, 0, 14.15914469713257, 41.99762985277868, 44.2901999575414, 18.30950700766612, 29.09245074905994, 8.905734062483303, 10.85018311856239, 12.21370809532799, 16.2240692932868, 21.8970263545853),
    @(13.63888885408329, 6.741981472312479, 0, 14.15983214881157, 41.9998179852342, 44.29466405332295, 18.31181948956764, 29.09641772138794, 8.905649375255704, 10.84548967039184, 12.2139764487257, 16.22336574005418, 21.89985725277736),
    @(13.68670388318942, 6.762859849183276, 0, 14.15519263269265, 41.9853782881293, 44.26460479662228, 18.29600069214426, 29.06929033062707, 8.906239774590206, 10.87822579546776, 12.21224994254083, 16.22846051976368, 21.88037016968154),
    @(13.90201254310768, 6.852391797846806, 0, 14.1379840371264, 41.94017794223334, 44.15492889315809, 18.23194569138387, 28.95967679355136, 8.90890818118242, 11.02639096298073, 12.20799314429841, 16.25614448968496, 21.79838862154872),
    @(14.33648648042805, 7.018249384931366, 0, 14.11591738266894, 41.91805803510381, 44.02277455841585, 18.12700245312835, 28.78104389187593, 8.914361520433165, 11.32818552661171, 12.21165269730326, 16.32883800731289, 21.65195040166476),
    @(14.66023504051899, 7.133869046546884, 0, 14.10677511974377, 41.94226396091089, 43.97650465775411, 18.06250245693345, 28.67197904308444, 8.91849748608302, 11.55486058470826, 12.22158771163982, 16.39318150160486, 21.55301147166065),
    @(14.80788411737194, 7.18504979414435, 0, 14.10413979799801, 41.9620232658357, 43.9665313381925, 18.0358963460822, 28.62718529938486, 8.920406065246848, 11.65864630740833, 12.22766316564535, 16.42475992652178, 21.50986153604633),
    @(14.86379798714963, 7.204221868821839, 0, 14.10335989096745, 41.97075873282508, 43.96434886279048, 18.02621466193757, 28.61091682946406, 8.921132591512935, 11.69800958110886, 12.23018582101411, 16.43704347113412, 21.49378760599577),
    @(14.85175663680673, 7.200102198267255, 0, 14.10351817813937, 41.96882175181179, 43.96474798027425, 18.02828228390118, 28.61438965999159, 8.920975954325911, 11.68952978841578, 12.22963267673293, 16.43438362020549, 21.4972376015156),
    @(14.81248442836178, 7.186631310850496, 0, 14.10407127108802, 41.96271685633216, 43.96631982873496, 18.03509194154945, 28.62583297048366, 8.920465761697354, 11.66188370358708, 12.22786626687141, 16.42576401606714, 21.50853379935572),
    @(14.78842786705162, 7.178352631102594, 0, 14.10443841865409, 41.95914045436292, 43.9674902745001, 18.03931429926887, 28.63293272364884, 8.92015374412858, 11.6449567198359, 12.22681314913264, 16.42052645007817, 21.51548765575681),
    @(14.65058893479961, 7.130495316119463, 0, 14.10697791303769, 41.94114827707666, 43.97737945835015, 18.06429628208826, 28.67500348605871, 8.918373292100958, 11.54808847058676, 12.22122183079744, 16.39116368179712, 21.55586870842333),
    @(14.56608813543617, 7.100769995401775, 0, 14.10892529380391, 41.93234802087266, 43.98628422178763, 18.08032259386059, 28.70204772509945, 8.917287907919462, 11.48880990803059, 12.21818906007516, 16.37373729181778, 21.58111620419401),
    @(14.51752375529871, 7.083539886038559, 0, 14.11018884453314, 41.92811009186512, 43.99244841269292, 18.08979795245146, 28.71805653097537, 8.91666618172983, 11.45477956406532, 12.21659125614995, 16.36393176284972, 21.59581283619784),
    @(14.50108884378981, 7.077683411970737, 0, 14.11064133319871, 41.92681681497356, 43.99471447028075, 18.09305035831992, 28.72355473225925, 8.916456118595717, 11.44326971551387, 12.21607549098302, 16.36064935515262, 21.60081894457516),
    @(14.57507980764295, 7.103948100353692, 0, 14.10870315077279, 41.93319960107645, 43.98522839437207, 18.07858991959939, 28.69912186023678, 8.917403184952056, 11.49511372391856, 12.21849674691782, 16.37556987842527, 21.57841046784414),
    @(14.82401996818105, 7.19059375086732, 0, 14.10390290557958, 41.96447604834033, 43.96581486462069, 18.03308110021722, 28.6224529532154, 8.920615515689255, 11.6700026223184, 12.22837909326045, 16.42828702224416, 21.50520862176004),
    @(14.98670752990844, 7.246001284793724, 0, 14.10203600256832, 41.99221781106734, 43.96241927333751, 18.00563197121582, 28.57639058713087, 8.922736991913869, 11.78464761786375, 12.23613095694816, 16.46463465869478, 21.4589169936685),
    @(14.89989565912956, 7.216542666170294, 0, 14.10291652841391, 41.97674534133976, 43.96338104907762, 18.0200721873862, 28.60060455582859, 8.9216027405246, 11.72343913194945, 12.23187591636566, 16.44506415702669, 21.48348225660882),
    @(14.57101461806087, 7.10251171674454, 0, 14.10880313310164, 41.93281204252032, 43.98570247979427, 18.07937244669836, 28.70044320868377, 8.917351061103121, 11.49226361132884, 12.21835718749282, 16.37474070163976, 21.57963316574463),
    @(14.21794058537319, 6.974453428110629, 0, 14.12064161338248, 41.91692706251407, 44.04961681976206, 18.15317945326508, 28.82547682173624, 8.912863152748345, 11.24553440396808, 12.20938441280656, 16.3072300791162, 21.6900408195713)
)

$startRow = 2
$startCol = 2  # Column B

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $row[$c]
    }
}

Write-Host "Updated loading_percent values for 380 kV case"